$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text: H1 (Affiliation Mismatch Note -> IsTrueIND)
# and J1 (SalesGroupName -> RepFirm)
$ws.Range("H1").Value = "IsTrueIND"
$ws.Range("J1").Value = "RepFirm"

# Give I1:J1 the same (bordered) formatting already used by the rest of the
# header row, by copying formats only from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Give I2:J5 the same (bordered) formatting already used by the rest of the
# data rows, by copying formats only from H2.
$ws.Range("H2").Copy()
$ws.Range("I2:J5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Reset the active selection to A1 so no stale selection is persisted.
$ws.Range("A1").Select()
